$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Applying functional directional light to drawn geometry." (row 30) as complete
$ws.Range("F30").Value = "X"

# Mark "Applying applicable color map texturing to drawn geometry." (row 18) as complete
$ws.Range("F18").Value = "X"

# Mark "Drawing indexed model loaded from file." (row 5) as complete
$ws.Range("F5").Value = "X"

# Mark "Demonstrates dynamic change in direction of directional lighting." (row 34)
# as achieved on Milestone I and complete
$ws.Range("E34").Value = "I"
$ws.Range("F34").Value = "X"

# Update the active selection to reflect where the student was working
# when these rows were completed
$ws.Range("A26").Select()

$wb.Save()
